# Add a new row of footer-data to the "pages_with_footer" worksheet and
# update the active selection to B12, as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_footer")

# Append the new record in row 8 (first sheet currently has data in A1:C7).
$ws.Range("A8").Value = "/news-events/press-releases/2018/oropharyngeal-hpv-cisplatin"
$ws.Range("B8").Value = "Article"
$ws.Range("C8").Value = "English"

# Reflect the updated selection that was captured when the workbook was saved.
$ws.Activate()
$ws.Range("B12").Select()
